# Apply updated crypto price/volume values as captured in the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.106.84"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "1.822.30"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -1.07%  "
$ws.Range("D5").Value = "'310.78"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "'0.4224"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").Value = "'0.3666"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("D9").Value = "'0.07228"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "'0.8479"
$ws.Range("E10").Value = "  -3.64%  "
$ws.Range("D11").Value = "'20.91"
$ws.Range("E11").Value = "  -3.39%  "
$ws.Range("D12").Value = "1.826.73"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "'6.673"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "'0.07090"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "'5.284"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "'89.46"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D18").Value = "'0.000008829"
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").Value = "27.270.86"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").Value = "'5.110"
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").Value = "'10.84"
$ws.Range("E23").Value = "  -2.65%  "
$ws.Range("D24").Value = "2.064.67"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").Value = "'152.02"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("D27").Value = "'2.204"
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("D28").Value = "'18.34"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "'5.211"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("D30").Value = "'116.27"
$ws.Range("E30").Value = "  -3.14%  "
$ws.Range("D31").Value = "'0.08813"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "'1.185"
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("D33").Value = "'0.7435"
$ws.Range("E33").Value = "  -4.50%  "
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'4.432"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "'1.100"
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").Value = "'0.01956"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "'0.05241"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").Value = "'7.257"
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").Value = "'2.875"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "'0.5026"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").Value = "'8.588"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "'10.58"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "'106.40"
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("D47").Value = "'0.4734"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "'1.659"
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("D51").Value = "'1.880"
$ws.Range("E51").Value = "  +0.52%  "
